$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 37 and 38 swap ranking order (Algorand moves up to rank 35,
# TheSandbox moves down to rank 36) plus updated price/volume figures.
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$d37 = $ws.Range("D37")
$d37.NumberFormat = "@"
$d37.Value = "0.2113"
$ws.Range("E37").Value = "  -2.74%  "

$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$d38 = $ws.Range("D38")
$d38.NumberFormat = "@"
$d38.Value = "0.6424"
$ws.Range("E38").Value = "  -2.18%  "

# Price / 1h-volume refresh for every other row in the table.
$ws.Range("D2").Value = "27.603.60"
$ws.Range("E2").Value = "  -2.31%  "

$ws.Range("D3").Value = "1.757.23"
$ws.Range("E3").Value = "  -3.10%  "

$d = $ws.Range("D4")
$d.NumberFormat = "@"
$d.Value = "1.003"
$ws.Range("E4").Value = "  +0.03%  "

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "324.79"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("E6").Value = "  +0.05%  "

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "0.4463"
$ws.Range("E7").Value = "  +2.92%  "

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = "0.3617"
$ws.Range("E8").Value = "  -1.14%  "

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = "0.07549"
$ws.Range("E9").Value = "  -1.66%  "

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "42.11"
$ws.Range("E10").Value = "  -6.14%  "

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "1.107"
$ws.Range("E11").Value = "  -3.00%  "

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "1.002"
$ws.Range("E12").Value = "  +0.04%  "

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "20.78"
$ws.Range("E13").Value = "  -5.58%  "

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "6.063"
$ws.Range("E14").Value = "  -3.62%  "

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = "7.221"
$ws.Range("E15").Value = "  -3.84%  "

$ws.Range("D16").Value = "1.756.40"
$ws.Range("E16").Value = "  -3.85%  "

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = "92.69"
$ws.Range("E17").Value = "  -1.12%  "

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = "0.00001066"
$ws.Range("E18").Value = "  -1.34%  "

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = "0.06418"
$ws.Range("E19").Value = "  -2.22%  "

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = "1.000"
$ws.Range("E20").Value = "  -0.06%  "

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = "17.08"
$ws.Range("E21").Value = "  -2.25%  "

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = "5.849"
$ws.Range("E22").Value = "  -6.41%  "

$ws.Range("D23").Value = "27.645.86"
$ws.Range("E23").Value = "  -2.25%  "

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = "11.26"
$ws.Range("E24").Value = "  -2.65%  "

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "2.101"
$ws.Range("E25").Value = "  +2.03%  "

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "162.58"
$ws.Range("E26").Value = "  -0.10%  "

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "20.43"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").Value = "1.957.66"
$ws.Range("E28").Value = "  -3.77%  "

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = "2.137"
$ws.Range("E29").Value = "  -6.31%  "

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = "125.94"
$ws.Range("E30").Value = "  -2.13%  "

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = "1.097"
$ws.Range("E31").Value = "  -9.20%  "

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = "3.660"
$ws.Range("E32").Value = "  +5.37%  "

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "5.572"
$ws.Range("E33").Value = "  -6.19%  "

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = "0.09016"
$ws.Range("E34").Value = "  -1.63%  "

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "12.17"
$ws.Range("E35").Value = "  -6.19%  "

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "0.02301"
$ws.Range("E36").Value = "  -2.07%  "

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = "0.06006"
$ws.Range("E39").Value = "  -3.03%  "

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = "4.944"
$ws.Range("E40").Value = "  -4.73%  "

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = "1.186"
$ws.Range("E41").Value = "  -0.55%  "

$ws.Range("E42").Value = "  +0.02%  "

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "1.400"
$ws.Range("E43").Value = "  -2.23%  "

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = "7.880"
$ws.Range("E44").Value = "  -2.82%  "

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "13.30"
$ws.Range("E45").Value = "  -4.32%  "

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "0.5926"
$ws.Range("E46").Value = "  -2.85%  "

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "3.714"
$ws.Range("E47").Value = "  -0.95%  "

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = "1.981"
$ws.Range("E48").Value = "  -1.66%  "

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = "121.95"
$ws.Range("E49").Value = "  -2.86%  "

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "1.168"
$ws.Range("E50").Value = "  +1.05%  "

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = "0.06875"
$ws.Range("E51").Value = "  -1.87%  "
